$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q3" immediately before the existing
#    "2022-Q2" sheet (right after "总计").
# ---------------------------------------------------------------------------
$existingQ2 = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($existingQ2)
$newSheet.Name = "2022-Q3"

# NOTE: sheet variables in this host resolve by tab position, so after the
# insert, the old `$existingQ2` handle now actually refers to the newly
# inserted sheet (it kept pointing at tab index 2). Re-fetch the real
# "2022-Q2" sheet by name so later reads/copies come from the right place.
$existingQ2 = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Item("2022-Q3")

# Copy the header row (values) from the "2022-Q2" sheet - every quarter
# sheet shares the same column headers. Column A has no header (the source
# sheet has no A1 cell at all), so only touch B1:H1.
$existingQ2.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Re-apply the header formatting (bold font + border) on top of the values.
$existingQ2.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the row-2/row-3 style used on every other quarter sheet: column A
# carries the bold bordered "index" style, columns B-H stay on the default
# style. Grab that style from the source sheet's A2 cell.
$existingQ2.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)
$newSheet.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Columns B, D, E, F, G hold digit-only strings ("001959", "4.62", "93.28",
# ...) that Excel's normal value-coercion would otherwise turn into real
# numbers (and, for the fund codes, silently drop the leading zero). Force
# those columns to Text first so the literal strings survive, matching the
# source data (fund name in column C is never numeric-looking, so it is
# left alone and keeps the default style).
$newSheet.Range("B2:B3").NumberFormat = "@"
$newSheet.Range("D2:G3").NumberFormat = "@"

# Row 2: 001959 / 华商乐享互联灵活配置混合A
$newSheet.Cells.Item(2, 1).Value = 0
$newSheet.Cells.Item(2, 2).Value = "001959"
$newSheet.Cells.Item(2, 3).Value = "华商乐享互联灵活配置混合A"
$newSheet.Cells.Item(2, 4).Value = "4.62"
$newSheet.Cells.Item(2, 5).Value = "93.28"
$newSheet.Cells.Item(2, 6).Value = "2.59"
$newSheet.Cells.Item(2, 7).Value = "0.1197"
$newSheet.Cells.Item(2, 8).Value = 9

# Row 3: 013142 / 华商乐享互联灵活配置混合C
$newSheet.Cells.Item(3, 1).Value = 1
$newSheet.Cells.Item(3, 2).Value = "013142"
$newSheet.Cells.Item(3, 3).Value = "华商乐享互联灵活配置混合C"
$newSheet.Cells.Item(3, 4).Value = "1.08"
$newSheet.Cells.Item(3, 5).Value = "93.28"
$newSheet.Cells.Item(3, 6).Value = "2.59"
$newSheet.Cells.Item(3, 7).Value = "0.0280"
$newSheet.Cells.Item(3, 8).Value = 9

# ---------------------------------------------------------------------------
# 2. Add the matching summary row to the "总计" sheet, right under the
#    header, pushing the other quarters down by one row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Preserve the "index" style (bold + border) on the row that is about to
# become the new last row (row 8) before we push values into it.
$total.Range("A7").Copy()
$total.Range("A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Shift the existing quarter rows down by one (bottom-up so we never
# clobber a value before reading it).
for ($r = 7; $r -ge 2; $r--) {
    $b = $total.Cells.Item($r, 2).Value2
    $c = $total.Cells.Item($r, 3).Value2
    $d = $total.Cells.Item($r, 4).Value2
    $dest = $r + 1
    $total.Cells.Item($dest, 1).Value = $r - 1
    $total.Cells.Item($dest, 2).Value = $b
    $total.Cells.Item($dest, 3).Value = $c
    $total.Cells.Item($dest, 4).Value = $d
}

# Write the new "2022-Q3" summary row.
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 2
$total.Cells.Item(2, 4).Value = 0.15
